$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "2024.11.01. 월간"
$ws.Range("B7").Value = "2024년 12월 08일 16시 47분 19초"
$ws.Range("C10").Value = "'96"
$ws.Range("D10").Value = "'41.38"
$ws.Range("C11").Value = "'136"
$ws.Range("D11").Value = "'58.62"
$ws.Range("C14").Value = "'1"
$ws.Range("D14").Value = "'0.43"
$ws.Range("C15").Value = "'3"
$ws.Range("D15").Value = "'1.29"
$ws.Range("C16").Value = "'31"
$ws.Range("D16").Value = "'13.36"
$ws.Range("C17").Value = "'33"
$ws.Range("D17").Value = "'14.22"
$ws.Range("C18").Value = "'15"
$ws.Range("D18").Value = "'6.47"
$ws.Range("C19").Value = "'19"
$ws.Range("D19").Value = "'8.19"
$ws.Range("C20").Value = "'7"
$ws.Range("D20").Value = "'3.02"
$ws.Range("C21").Value = "'22"
$ws.Range("D21").Value = "'9.48"
$ws.Range("C22").Value = "'8"
$ws.Range("D22").Value = "'3.45"
$ws.Range("C23").Value = "'17"
$ws.Range("D23").Value = "'7.33"
$ws.Range("C24").Value = "'8"
$ws.Range("D24").Value = "'3.45"
$ws.Range("C25").Value = "'14"
$ws.Range("D25").Value = "'6.03"
$ws.Range("C26").Value = "'19"
$ws.Range("D26").Value = "'8.19"
$ws.Range("C27").Value = "'11"
$ws.Range("D27").Value = "'4.74"
$ws.Range("C28").Value = "'5"
$ws.Range("D28").Value = "'2.16"
$ws.Range("C29").Value = "'12"
$ws.Range("D29").Value = "'5.17"
$ws.Range("C30").Value = "'1"
$ws.Range("D30").Value = "'0.43"
$ws.Range("C31").Value = "'4"
$ws.Range("D31").Value = "'1.72"
$ws.Range("C32").Value = "'1"
$ws.Range("D32").Value = "'0.43"
$ws.Range("C33").Value = "'1"
$ws.Range("D33").Value = "'0.43"
